$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "numero" (A) and "descricao" (B) values for rows 4 through 25,
# keeping the existing style but leaving the cells blank.
$ws.Range("A4:B25").ClearContents()
